$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "Package info" (sheet1): update package metadata, add a new
# "Name" row, and narrow column B.
# -----------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("Package info")

$wsInfo.Cells.Item(3,2).Value = "AGGREGATE"
$wsInfo.Cells.Item(4,2).Value = "1.2.1"
$wsInfo.Cells.Item(5,2).Value = "2.35.6"
$wsInfo.Cells.Item(6,1).Value = "DHIS2 build"
$wsInfo.Cells.Item(6,2).Value = "9979080"
$wsInfo.Cells.Item(7,1).Value = "Last updated"
$wsInfo.Cells.Item(7,2).Value = "20210915T115527"

# New row 8 ("Name") - copy banding/format from row 6 (same parity) first
$wsInfo.Range("A6:B6").Copy()
$wsInfo.Range("A8:B8").PasteSpecial(-4122)
$wsInfo.Cells.Item(8,1).Value = "Name"
$wsInfo.Cells.Item(8,2).Value = "COVIDVAC_AGGREGATE_V1.2.1_2.35.6-en"

$wsInfo.Columns.Item(2).ColumnWidth = 36.833333333333336
